$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typos in stuurvariabelen descriptions (Sinuositeitsklasse / natprofiel gemaaid)
$le = [char]0x2264
$d9text = "1 = Recht (sinuositeit=1); 2 = Gestrekt (1<sinuositeit" + $le + "1,05);  3 = Gestrekt (1.05<sinuositeit" + $le + "1,25); 4 = Gestrekt (1.25<sinuositeit" + $le + "1,5); 5 = Gestrekt (sinuositeit>1,5)"
$ws.Range("D9").Value = $d9text
$ws.Range("D14").Value = "'(- Percentage nat profiel gemaaid / 200)"

# Match the final selection state left behind in the saved workbook
$ws.Range("D24").Select() | Out-Null
